# Update column G ("K") values on Sheet1 rows 2-34 as per regenerated save_data
# (commit: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 4
    3  = 6
    4  = 9
    5  = 5
    6  = 6
    7  = 5
    8  = 7
    9  = 3
    10 = 5
    11 = 7
    12 = 5
    13 = 7
    14 = 4
    15 = 2
    16 = 5
    17 = 8
    18 = 7
    19 = 6
    20 = 3
    21 = 3
    22 = 8
    23 = 3
    24 = 8
    25 = 7
    26 = 7
    27 = 5
    28 = 10
    29 = 7
    30 = 11
    31 = 6
    32 = 4
    33 = 5
    34 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
